$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-88 down to 67-89.
$ws.Rows.Item(66).Insert()

# Fill in the new row 66 with the new weekly record.
$ws.Range("A66").Value = 4
$ws.Range("B66").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C66").Value = "Los Lagos"
$ws.Range("D66").Value = 44524
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = 100112022
$ws.Range("G66").Value = "Arveja Verde"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 15
$ws.Range("K66").Value = 18000
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = 18000
$ws.Range("N66").Value = "`$/saco 25 kilos"
$ws.Range("O66").Value = "Región del Maule"
$ws.Range("P66").Value = 720
$ws.Range("Q66").Value = 25
$ws.Range("R66").Value = "Hortaliza"
